# project-status.xlsx update
# Task "add admin user" (routes row 4) is done; route/controller for /api/user/add
# was written and tested. New current task is "add faculty" (routes row 5).
# Logged the completed work as a new row in the "tasks" sheet.

$wb = $excel.ActiveWorkbook

$routes = $wb.Worksheets.Item("routes")
$tasks  = $wb.Worksheets.Item("tasks")

# --- tasks sheet: log the completed /api/user/add work ---
$tasks.Range("A5").Value = 45555
$tasks.Range("A5").NumberFormat = "d-mmm"
$tasks.Range("B5").Value = "Jasdeep"
$tasks.Range("C5").Value = "Wrote and tested /api/user/add route to add admin user"
$tasks.Range("D5").Value = "test cases for faculty have to be included"

# --- routes sheet: row 4 ("add admin user") moves from "current task" to "done" ---
$routes.Range("F4").Value = "done"
$routes.Range("G4").Value = 45555
$routes.Range("H4").Value = "admin or super admin can create a new admin"
$routes.Range("J4").Value = "user-add.test.js"

# --- routes sheet: row 5 ("add faculty") becomes the new "current task" ---
$routes.Range("E5").Value = "Jasdeep"
$routes.Range("F5").Value = "current task"
$routes.Range("G5").Value = 45555
$routes.Range("G5").NumberFormat = "d-mmm"

# --- selection / active sheet bookkeeping (user ended up on "routes", cell D5) ---
$tasks.Range("D5").Select() | Out-Null
$routes.Activate() | Out-Null
$routes.Range("D5").Select() | Out-Null

Write-Host "edit applied"
